$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cell M102 (timestamp precision correction) ---
$ws.Range("M102").Value = 45905.77089622685
$ws.Range("M102").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Append new rows 103-112 (newly logged experiment results) ---

# Row 103
$ws.Range("A103").Value = "Fucntionality_test_AIDS_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B103").Value = "AIDS"
$ws.Range("C103").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D103").Value = 0.2
$ws.Range("E103").Value = 0.0000006855555555555555
$ws.Range("F103").Value = 0.0000001162037037037037
$ws.Range("G103").Value = 0.9925
$ws.Range("H103").Value = 0.9924467959850469
$ws.Range("I103").Value = 0.9925698757763974
$ws.Range("J103").Value = 0.9925
$ws.Range("K103").Value = 0.9814814814814815
$ws.Range("L103").Value = 45910.64188774306
$ws.Range("L103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M103").Value = 45910.64188697917
$ws.Range("M103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N103").Value = "GEDLIB_Calculator"
$ws.Range("O103").Value = "Simple Train-Test Split"

# Row 104
$ws.Range("A104").Value = "Fucntionality_test_AIDS_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B104").Value = "AIDS"
$ws.Range("C104").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D104").Value = 0.2
$ws.Range("E104").Value = 0.0000006855555555555555
$ws.Range("F104").Value = 0.0000001162037037037037
$ws.Range("G104").Value = 0.995
$ws.Range("H104").Value = 0.9949765625000001
$ws.Range("I104").Value = 0.9950311526479751
$ws.Range("J104").Value = 0.995
$ws.Range("K104").Value = 0.9876543209876543
$ws.Range("L104").Value = 45910.64188774306
$ws.Range("L104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M104").Value = 45910.64188697917
$ws.Range("M104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N104").Value = "GEDLIB_Calculator"
$ws.Range("O104").Value = "Hyperparameter Tuning (grid)"

# Row 105
$ws.Range("A105").Value = "Fucntionality_test_PTC_FR_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B105").Value = "PTC_FR"
$ws.Range("C105").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D105").Value = 0.2
$ws.Range("E105").Value = 0.0000001629976851851852
$ws.Range("F105").Value = 0.00000001989583333333333
$ws.Range("G105").Value = 0.5070422535211268
$ws.Range("H105").Value = 0.5082502365633458
$ws.Range("I105").Value = 0.5096530401923738
$ws.Range("J105").Value = 0.5070422535211268
$ws.Range("K105").Value = 0.4926108374384236
$ws.Range("L105").Value = 45910.6429222338
$ws.Range("L105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M105").Value = 45910.64292207176
$ws.Range("M105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N105").Value = "GEDLIB_Calculator"
$ws.Range("O105").Value = "Simple Train-Test Split"

# Row 106
$ws.Range("A106").Value = "Fucntionality_test_PTC_FR_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B106").Value = "PTC_FR"
$ws.Range("C106").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D106").Value = 0.2
$ws.Range("E106").Value = 0.0000001629976851851852
$ws.Range("F106").Value = 0.00000001989583333333333
$ws.Range("G106").Value = 0.5915492957746479
$ws.Range("H106").Value = 0.4397357596908887
$ws.Range("I106").Value = 0.3499305693314819
$ws.Range("J106").Value = 0.5915492957746479
$ws.Range("K106").Value = 0.5
$ws.Range("L106").Value = 45910.6429222338
$ws.Range("L106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M106").Value = 45910.64292207176
$ws.Range("M106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N106").Value = "GEDLIB_Calculator"
$ws.Range("O106").Value = "Hyperparameter Tuning (grid)"

# Row 107
$ws.Range("A107").Value = "Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B107").Value = "MUTAG"
$ws.Range("C107").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D107").Value = 0.2
$ws.Range("E107").Value = 0.0000001471296296296296
$ws.Range("F107").Value = 0.000000007083333333333334
$ws.Range("G107").Value = 0.7894736842105263
$ws.Range("H107").Value = 0.7989203778677463
$ws.Range("I107").Value = 0.8258145363408521
$ws.Range("J107").Value = 0.7894736842105263
$ws.Range("K107").Value = 0.7928571428571429
$ws.Range("L107").Value = 45910.64468162037
$ws.Range("L107").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M107").Value = 45910.64468146991
$ws.Range("M107").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N107").Value = "GEDLIB_Calculator"
$ws.Range("O107").Value = "Simple Train-Test Split"

# Row 108
$ws.Range("A108").Value = "Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B108").Value = "MUTAG"
$ws.Range("C108").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D108").Value = 0.2
$ws.Range("E108").Value = 0.0000001471296296296296
$ws.Range("F108").Value = 0.000000007083333333333334
$ws.Range("G108").Value = 0.6842105263157895
$ws.Range("H108").Value = 0.7030075187969925
$ws.Range("I108").Value = 0.7801169590643275
$ws.Range("J108").Value = 0.6842105263157895
$ws.Range("K108").Value = 0.7214285714285714
$ws.Range("L108").Value = 45910.64468162037
$ws.Range("L108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M108").Value = 45910.64468146991
$ws.Range("M108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N108").Value = "GEDLIB_Calculator"
$ws.Range("O108").Value = "Hyperparameter Tuning (grid)"

# Row 109
$ws.Range("A109").Value = "Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B109").Value = "MUTAG"
$ws.Range("C109").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D109").Value = 0.2
$ws.Range("E109").Value = 0.0000002213078703703704
$ws.Range("F109").Value = 0.000000006261574074074074
$ws.Range("G109").Value = 0.7105263157894737
$ws.Range("H109").Value = 0.7119788719109601
$ws.Range("I109").Value = 0.7144138755980861
$ws.Range("J109").Value = 0.7105263157894737
$ws.Range("K109").Value = 0.7028985507246377
$ws.Range("L109").Value = 45910.64681074074
$ws.Range("L109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M109").Value = 45910.64681050926
$ws.Range("M109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N109").Value = "GEDLIB_Calculator"
$ws.Range("O109").Value = "Simple Train-Test Split"

# Row 110
$ws.Range("A110").Value = "Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B110").Value = "MUTAG"
$ws.Range("C110").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D110").Value = 0.2
$ws.Range("E110").Value = 0.0000002213078703703704
$ws.Range("F110").Value = 0.000000006261574074074074
$ws.Range("G110").Value = 0.8157894736842105
$ws.Range("H110").Value = 0.8111946532999164
$ws.Range("I110").Value = 0.8178137651821862
$ws.Range("J110").Value = 0.8157894736842105
$ws.Range("K110").Value = 0.7898550724637682
$ws.Range("L110").Value = 45910.64681074074
$ws.Range("L110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M110").Value = 45910.64681050926
$ws.Range("M110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N110").Value = "GEDLIB_Calculator"
$ws.Range("O110").Value = "Hyperparameter Tuning (grid)"

# Row 111
$ws.Range("A111").Value = "Fucntionality_test_Letter-high_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B111").Value = "Letter-high"
$ws.Range("C111").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D111").Value = 0.2
$ws.Range("E111").Value = 0.00000555255787037037
$ws.Range("F111").Value = 0.0000001163310185185185
$ws.Range("G111").Value = 0.8755555555555555
$ws.Range("H111").Value = 0.9075116598079562
$ws.Range("I111").Value = 0.9628330876493202
$ws.Range("J111").Value = 0.8755555555555555
$ws.Range("K111").Value = 0.8847234094517036
$ws.Range("L111").Value = 45910.65059689815
$ws.Range("L111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M111").Value = 45910.65059133102
$ws.Range("M111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N111").Value = "GEDLIB_Calculator"
$ws.Range("O111").Value = "Simple Train-Test Split"

# Row 112
$ws.Range("A112").Value = "Fucntionality_test_Letter-high_with_SVC_Simple_Prototype_GED_poly"
$ws.Range("B112").Value = "Letter-high"
$ws.Range("C112").Value = "SVC_Simple_Prototype_GED_poly"
$ws.Range("D112").Value = 0.2
$ws.Range("E112").Value = 0.00000555255787037037
$ws.Range("E112").NumberFormat = "0"
$ws.Range("F112").Value = 0.0000001163310185185185
$ws.Range("F112").NumberFormat = "0"
$ws.Range("G112").Value = 0.9577777777777777
$ws.Range("H112").Value = 0.9371219573716736
$ws.Range("I112").Value = 0.9173382716049383
$ws.Range("J112").Value = 0.9577777777777777
$ws.Range("K112").Value = 0.5
$ws.Range("L112").Value = 45910.65059689566
$ws.Range("L112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M112").Value = 45910.65059132913
$ws.Range("M112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N112").Value = "GEDLIB_Calculator"
$ws.Range("O112").Value = "Hyperparameter Tuning (grid)"
